function Replace-Text {
    param($old, $new)
    $d = $word.ActiveDocument
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Text = $new
    } else {
        Write-Output "NOT FOUND: $old"
    }
}

Replace-Text 'English' 'Inglese'
Replace-Text ' / Portuguese / French / Thai / Vietnamese / Spanish' ' / Portoghese / Francese / Thai / Vietnamita / Spagnolo'
Replace-Text 'English' 'Inglese'
Replace-Text 'Brief' 'Breve'
Replace-Text 'An email sent to partners in the target country who have sent their documents for review. It will be sent via customer.io' 'Un''email inviata ai partner nel paese target che hanno inviato i loro documenti per la revisione. Sarà inviata tramite customer.io'
Replace-Text 'Target audience' 'Pubblico target'
Replace-Text 'Invited partners who have submitted their documents' 'Partner invitati che hanno presentato i loro documenti'
Replace-Text 'Subject line' 'Oggetto'
Replace-Text ' — we got your docs!  ' ' — abbiamo ricevuto i vostri documenti!  '
Replace-Text 'Thank you for submitting your documents' 'Grazie per aver inviato i documenti'
Replace-Text 'Hi ' 'Ciao '
Replace-Text 'Thank you for providing us with your documents for the upcoming ' 'Grazie per aver fornito i documenti per il prossimo '
Replace-Text '. Based on the information you’ve given us, we’ll make the necessary arrangements, including accommodation and transportation.' '. Sulla base delle informazioni ricevute, ci occuperemo dell''organizzazione necessaria, che comprenderà alloggio e trasporto.'
Replace-Text 'We’re currently reviewing your documents and will reach out to you if we need anything else. ' 'Stiamo esaminando i tuoi documenti e ti contatteremo se avremo bisogno di altro. '
Replace-Text 'If you have any questions, please contact us via ' 'Se hai domande, non esitare a contattarci tramite '
Replace-Text ' or ' ' o '
Replace-Text 'If you have any questions, please contact your country manager, ' 'In caso di domande, contatta il tuo country manager, '
Replace-Text ', at ' ', all''indirizzo '
Replace-Text ' or ' ' o al numero '
Replace-Text 'We look forward to seeing you at ' 'Non vediamo l''ora di incontrarti all''evento '

# Comment text update (comments story isn't reachable via Find, so set Range.Text directly)
$d = $word.ActiveDocument
$d.Comments.Item(1).Range.Text = 'scegli uno qualsiasi'

